$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 2 values that changed
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1

# Delete rows 3 to 5 (data no longer present), which also removes the
# "Neutrophils" shared string since it's no longer referenced anywhere.
$ws.Rows("3:5").Delete()
